$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 959.625
$ws.Range("I17").Value = 767.4
$ws.Range("J17").Value = 1280
$ws.Range("K17").Value = 2302.2
$ws.Range("L17").Value = 3840
$ws.Range("M17").Value = -2134.2
$ws.Range("N17").Value = -4176
$ws.Range("H58").Value = 850.9231
$ws.Range("I58").Value = 415
$ws.Range("J58").Value = 3248.5
$ws.Range("K58").Value = 1245
$ws.Range("L58").Value = 9745.5
$ws.Range("M58").Value = -1095
$ws.Range("N58").Value = -10045.5
$ws.Range("H62").Value = 1902.6666
$ws.Range("J62").Value = 1854
$ws.Range("L62").Value = 1854
$ws.Range("N62").Value = -3102
$ws.Range("H65").Value = 1902.6666
$ws.Range("J65").Value = 1854
$ws.Range("L65").Value = 9270
$ws.Range("N65").Value = -15510
$ws.Range("H74").Value = 6090.1035
$ws.Range("I74").Value = 5319
$ws.Range("J74").Value = 16500
$ws.Range("K74").Value = 5319
$ws.Range("L74").Value = 16500
$ws.Range("M74").Value = -4383
$ws.Range("N74").Value = -18372
$ws.Range("H77").Value = 6090.1035
$ws.Range("I77").Value = 5319
$ws.Range("J77").Value = 16500
$ws.Range("K77").Value = 26595
$ws.Range("L77").Value = 82500
$ws.Range("M77").Value = -21915
$ws.Range("N77").Value = -91860
$ws.Range("H96").Value = 2209.4
$ws.Range("I96").Value = 538
$ws.Range("K96").Value = 1614
$ws.Range("M96").Value = -241
$ws.Range("H98").Value = 1389.2307
$ws.Range("I98").Value = 1389.2307
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1389.2307
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 108.7692999999999
$ws.Range("N98").ClearContents()
$ws.Range("H106").Value = 15218.091
$ws.Range("I106").Value = 2399.5
$ws.Range("K106").Value = 2399.5
$ws.Range("M106").Value = -1768.5
$ws.Range("H111").Value = 1042.6666
$ws.Range("I111").Value = 1042.6666
$ws.Range("K111").Value = 3127.9998
$ws.Range("M111").Value = -60.99980000000005
$ws.Range("H122").Value = 1389.2307
$ws.Range("I122").Value = 1389.2307
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4167.6921
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1717.6921
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 9644.799999999999
$ws.Range("I132").Value = 9644.799999999999
$ws.Range("K132").Value = 28934.4
$ws.Range("M132").Value = -26404.4
$ws.Range("H141").Value = 1051.76
$ws.Range("I141").Value = 1051.76
$ws.Range("K141").Value = 3155.28
$ws.Range("M141").Value = 2024.72

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 587555.5
$ws.Range("I32").Value = 618753.8
$ws.Range("K32").Value = 618753.8
$ws.Range("M32").Value = -618466.8
$ws.Range("H61").Value = 8335219
$ws.Range("I61").Value = 2128.75
$ws.Range("J61").Value = 25001400
$ws.Range("K61").Value = 2128.75
$ws.Range("L61").Value = 25001400
$ws.Range("M61").Value = -1916.75
$ws.Range("N61").Value = -25001824
$ws.Range("H98").Value = 21955
$ws.Range("J98").Value = 21955
$ws.Range("L98").Value = 21955
$ws.Range("N98").Value = -27945
$ws.Range("H102").Value = 3554.3
$ws.Range("I102").Value = 2886.25
$ws.Range("K102").Value = 2886.25
$ws.Range("M102").Value = -1264.25
$ws.Range("H122").Value = 2518.5
$ws.Range("I122").Value = 2531.8
$ws.Range("J122").Value = 2505.2
$ws.Range("K122").Value = 7595.400000000001
$ws.Range("L122").Value = 7515.599999999999
$ws.Range("M122").Value = -5145.400000000001
$ws.Range("N122").Value = -12415.6
$ws.Range("H132").Value = 2708.587
$ws.Range("I132").Value = 1699.8518
$ws.Range("K132").Value = 5099.555399999999
$ws.Range("M132").Value = -2569.555399999999
$ws.Range("H136").Value = 8335219
$ws.Range("I136").Value = 2128.75
$ws.Range("J136").Value = 25001400
$ws.Range("K136").Value = 6386.25
$ws.Range("L136").Value = 75004200
$ws.Range("M136").Value = -3836.25
$ws.Range("N136").Value = -75009300

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 38917.82
$ws.Range("I20").Value = 88144
$ws.Range("K20").Value = 88144
$ws.Range("M20").Value = -87897

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 75.86207
$ws.Range("I7").Value = 44.8
$ws.Range("J7").Value = 109.14286
$ws.Range("K7").Value = 44.8
$ws.Range("L7").Value = 109.14286
$ws.Range("M7").Value = 68.2
$ws.Range("N7").Value = -335.14286
$ws.Range("H31").Value = 1018242.75
$ws.Range("I31").Value = 1304029.8
$ws.Range("K31").Value = 1304029.8
$ws.Range("M31").Value = -1303734.8
$ws.Range("H34").Value = 1018242.75
$ws.Range("I34").Value = 1304029.8
$ws.Range("K34").Value = 1304029.8
$ws.Range("M34").Value = -1303827.8
$ws.Range("H132").Value = 4863.6294
$ws.Range("I132").Value = 3728.762
$ws.Range("J132").Value = 8835.666999999999
$ws.Range("K132").Value = 11186.286
$ws.Range("L132").Value = 26507.001
$ws.Range("M132").Value = -8656.286
$ws.Range("N132").Value = -31567.001
$ws.Range("H134").Value = 2488.2812
$ws.Range("I134").Value = 2320.8667
$ws.Range("K134").Value = 6962.6001
$ws.Range("M134").Value = -4427.6001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1512.8
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8246
$ws.Range("H84").Value = 1512.8
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 18000
$ws.Range("N84").Value = -29232
$ws.Range("H92").Value = 841.9259
$ws.Range("I92").Value = 499.94736
$ws.Range("J92").Value = 1654.125
$ws.Range("K92").Value = 1499.84208
$ws.Range("L92").Value = 4962.375
$ws.Range("M92").Value = -251.8420799999999
$ws.Range("N92").Value = -7458.375
$ws.Range("H109").Value = 1333.9286
$ws.Range("J109").Value = 4000
$ws.Range("L109").Value = 12000
$ws.Range("N109").Value = -14080

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 58900
$ws.Range("J93").Value = 58900
$ws.Range("L93").Value = 58900
$ws.Range("N93").Value = -62644
$ws.Range("H132").Value = 16235.73
$ws.Range("I132").Value = 18114.348
$ws.Range("J132").Value = 1833
$ws.Range("K132").Value = 54343.04400000001
$ws.Range("L132").Value = 5499
$ws.Range("M132").Value = -51813.04400000001
$ws.Range("N132").Value = -10559

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4321.0625
$ws.Range("I7").Value = 4374.7
$ws.Range("K7").Value = 4374.7
$ws.Range("M7").Value = -4262.7
$ws.Range("H100").Value = 2975
$ws.Range("I100").Value = 2975
$ws.Range("K100").Value = 2975
$ws.Range("M100").Value = -2434
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H126").Value = 4321.0625
$ws.Range("I126").Value = 4374.7
$ws.Range("K126").Value = 13124.1
$ws.Range("M126").Value = -10654.1
$ws.Range("H136").Value = 6412820
$ws.Range("I136").Value = 3679058.2
$ws.Range("J136").Value = 25002398
$ws.Range("K136").Value = 11037174.6
$ws.Range("L136").Value = 75007194
$ws.Range("M136").Value = -11034624.6
$ws.Range("N136").Value = -75012294

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 21158.5
$ws.Range("I62").Value = 20730.25
$ws.Range("J62").Value = 21729.5
$ws.Range("K62").Value = 20730.25
$ws.Range("L62").Value = 21729.5
$ws.Range("M62").Value = -20106.25
$ws.Range("N62").Value = -22977.5
$ws.Range("H65").Value = 21158.5
$ws.Range("I65").Value = 20730.25
$ws.Range("J65").Value = 21729.5
$ws.Range("K65").Value = 103651.25
$ws.Range("L65").Value = 108647.5
$ws.Range("M65").Value = -100531.25
$ws.Range("N65").Value = -114887.5
$ws.Range("H126").Value = 3124.1765
$ws.Range("I126").Value = 3775
$ws.Range("K126").Value = 11325
$ws.Range("M126").Value = -8855
$ws.Range("H130").Value = 44980
$ws.Range("J130").Value = 44980
$ws.Range("L130").Value = 44980
$ws.Range("N130").Value = -55020
